$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-apply the percentage formula on Hoja1 as one range-level write so the
#    engine collapses the identical per-row formulas into a single
#    shared-formula group (t="shared"), matching a genuine fill-down.
# ---------------------------------------------------------------------------
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("C3:C15").Formula = '=B3/B$1*100'

# ---------------------------------------------------------------------------
# 2. Add the two new sheets, right after "Hoja1", in final tab order:
#    Hoja1, restaurant, retaurant_proudct
# ---------------------------------------------------------------------------
$restaurant = $wb.Worksheets.Add($null, $hoja1)
$restaurant.Name = "restaurant"

$retProd = $wb.Worksheets.Add($null, $restaurant)
$retProd.Name = "retaurant_proudct"

# ---------------------------------------------------------------------------
# 3. "restaurant" sheet content. Labels are entered in the same order the
#    original workbook's shared-string table shows them being introduced
#    (alto, acho, bottom card, imagen sup, radio) so newly-created shared
#    strings line up the same way.
# ---------------------------------------------------------------------------
$restaurant.Range("A1").Value = "alto"
$restaurant.Range("B1").Value = 438

$restaurant.Range("A2").Value = "acho"
$restaurant.Range("B2").Value = 202

$restaurant.Range("A5").Value = "bottom card"
$restaurant.Range("B5").Value = 295

$restaurant.Range("A4").Value = "imagen sup"
$restaurant.Range("B4").Value = 149

$restaurant.Range("A6").Value = "radio"
$restaurant.Range("B6").Value = 10

$restaurant.Range("C4").Formula = '=B4/B$1*100'
$restaurant.Range("D4").Formula = '=B1-B4'

$restaurant.Range("C5").Formula = '=B5/B$1*100'
$restaurant.Range("D5").Formula = '=B1*C5/100'

$restaurant.Range("C6").Formula = '=B6/B$1*100'

$restaurant.Range("B7").Value = 72.88
$restaurant.Range("C7").Formula = '=B7/B$1*100'

$restaurant.Range("B8").Value = 29
$restaurant.Range("C8").Formula = '=B8/B$1*100'

$restaurant.Range("B9").Value = 62
$restaurant.Range("C9").Formula = '=B9/B$1*100'

# Column C carries a 2-decimal number format; columns A and C keep their
# authored widths.
$restaurant.Range("C4:C9").NumberFormat = "0.00"
$restaurant.Columns.Item(1).ColumnWidth = 13.71
$restaurant.Columns.Item(3).ColumnWidth = 11.43

$restaurant.Activate()
[void]$restaurant.Range("C9").Select()

# ---------------------------------------------------------------------------
# 4. "retaurant_proudct" sheet content - same idea: labels entered in the
#    order the shared-string table introduces them (alt, alt card, alt add,
#    slider w, anc).
# ---------------------------------------------------------------------------
$retProd.Range("A1").Value = "alt"
$retProd.Range("B1").Value = 560
$retProd.Range("E1").Value = 559

$retProd.Range("A3").Value = "alt card"
$retProd.Range("B3").Value = 314
$retProd.Range("E3").Value = 42

$retProd.Range("A4").Value = "alt add "
$retProd.Range("B4").Value = 38

$retProd.Range("A5").Value = "slider w"
$retProd.Range("B5").Value = 144

$retProd.Range("A2").Value = "anc"
$retProd.Range("B2").Value = 264
$retProd.Range("E2").Value = 259

$retProd.Range("C3").Formula = '=B3*100/B$1'
$retProd.Range("F3").Formula = '=E3/E2*100'

$retProd.Range("C4").Formula = '=B4*100/B$1'

$retProd.Range("C5").Formula = '=B5*100/B2'

$retProd.Activate()
[void]$retProd.Range("F4").Select()

# ---------------------------------------------------------------------------
# 5. Window / tab state: "retaurant_proudct" ends up the active tab, and the
#    tab strip is scrolled so "restaurant" is the first visible sheet tab.
# ---------------------------------------------------------------------------
$retProd.Activate()
[void]$wb.Windows.Item(1).ScrollWorkbookTabs(1)
